$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 152, shifting existing rows 152..253 down to 153..254.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new record.
$ws.Cells.Item(152, 1).Value = 10
$ws.Cells.Item(152, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value = "La Araucanía"
$ws.Cells.Item(152, 4).Value = 44574
$ws.Cells.Item(152, 5).Value = 9
$ws.Cells.Item(152, 6).Value = 100112009
$ws.Cells.Item(152, 7).Value = "Acelga"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 205
$ws.Cells.Item(152, 11).Value = 7000
$ws.Cells.Item(152, 12).Value = 8000
$ws.Cells.Item(152, 13).Value = 7610
$ws.Cells.Item(152, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(152, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(152, 16).Value = 634
$ws.Cells.Item(152, 17).Value = 12
$ws.Cells.Item(152, 18).Value = "Hortaliza"
